$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -3
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -1
$ws.Range("F10").Value = 1
$ws.Range("F12").Value = -2
$ws.Range("F14").Value = -7
$ws.Range("F15").Value = -1
